# Actualización automática 2025-09-12 16:45:09
$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("M15").Value = 1451.52

# --- Sheet: VENTA MENSUAL ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F15").Value = 1451.52
$wsMensual.Range("F31").Value = 6242.92

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumpl.Range("D12").Value = 1451.52
$wsCumpl.Range("E12").Value = 16223.8286842162
$wsCumpl.Range("F12").Value = 0.08212115222915993

$wsCumpl.Range("D15").Value = 6317.210000000001
$wsCumpl.Range("E15").Value = 25390.54990313501
$wsCumpl.Range("F15").Value = 0.1992323021020292
